$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26; existing rows 26-62 shift down to 27-63.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with its data (matches the surrounding
# records for this market/category, with its own date/volume/price figures).
$ws.Cells.Item(26, 1).Value = 9
$ws.Cells.Item(26, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(26, 3).Value = "Metropolitana"
$ws.Cells.Item(26, 4).Value = 44799
$ws.Cells.Item(26, 5).Value = 13
$ws.Cells.Item(26, 6).Value = 100112035
$ws.Cells.Item(26, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 27
$ws.Cells.Item(26, 11).Value = 20000
$ws.Cells.Item(26, 12).Value = 20000
$ws.Cells.Item(26, 13).Value = 20000
$ws.Cells.Item(26, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(26, 15).Value = "Hijuelas"
$ws.Cells.Item(26, 16).Value = 1333
$ws.Cells.Item(26, 17).Value = 15
$ws.Cells.Item(26, 18).Value = "Hortaliza"
